# Final Checklist.xlsx edit
# Commit: "Add my hours to Final Checklist."
# (The other parts of the commit - deleting PatientID from the Appointment
#  table in the Dev Doc / Project Doc, and renaming Project Doc - touch
#  separate files outside this workbook and have no effect here.)
#
# Fill in the "Actual Hours" column (F) on Sheet1 for the checklist items
# that have been completed so far.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value = 2      # A development documentation with Excel
$ws.Range("F4").Value = 0.5    # An ERD
$ws.Range("F5").Value = 2      # A SQL database script
$ws.Range("F8").Value = 2      # A backup of the database with sample data
$ws.Range("F12").Value = 2     # A formal project document with Word
$ws.Range("F13").Value = 1     # Lessons Learned Document

# Move the active selection to where the author last clicked.
$ws.Range("E7").Select() | Out-Null
